$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# The workbook originally had 3 groups of 4 metric columns per model
# (AUC Lower, AUC Upper, Acc. Lower, Acc. Upper) in B:E, F:I, J:M.
# The fix adds a new "AUC" column (between AUC Lower/AUC Upper) and a new
# "Acc." column (between Acc. Lower/Acc. Upper) to each group, turning each
# 4-column group into a 6-column group (B:G, H:M, N:S), and appends a new
# "n" column (sample count per semester) at U.
# ---------------------------------------------------------------------------

# Insert the 6 new (blank) columns. Working right-to-left so that earlier
# column letters remain valid reference points as later inserts shift
# everything after them to the right.
$insertCols = @("M", "K", "I", "G", "E", "C")
foreach ($c in $insertCols) {
    $ws.Columns($c).Insert()
}

# Each inserted column takes on the column width of its immediate left
# neighbor (the same behaviour Excel itself shows for a manual column
# insert), instead of the sheet's default column width.
$ws.Columns("C").ColumnWidth = $ws.Columns("B").ColumnWidth
$ws.Columns("F").ColumnWidth = $ws.Columns("E").ColumnWidth
$ws.Columns("I").ColumnWidth = $ws.Columns("H").ColumnWidth
$ws.Columns("L").ColumnWidth = $ws.Columns("K").ColumnWidth
$ws.Columns("O").ColumnWidth = $ws.Columns("N").ColumnWidth
$ws.Columns("R").ColumnWidth = $ws.Columns("Q").ColumnWidth

# Row 2 sub-headers for the newly inserted columns.
$ws.Range("C2").Value = "AUC"
$ws.Range("F2").Value = "Acc."
$ws.Range("I2").Value = "AUC"
$ws.Range("L2").Value = "Acc."
$ws.Range("O2").Value = "AUC"
$ws.Range("R2").Value = "Acc."

# New data values (plain AUC / plain Accuracy, computed from probability
# rather than the rounded prediction) for each of the three model groups.
$newData = @{
    3  = @(0.83684895122705505, 0.79100145137880895, 0.86335025196019499, 0.79390420899854797, 0.49795339191180399, 0.63570391872278598)
    4  = @(0.87889225682006,    0.82165605095541405, 0.88236069344770895, 0.81687898089171895, 0.49790935339715098, 0.69585987261146498)
    5  = @(0.84855906750571997, 0.83539823008849501, 0.86037614416475905, 0.82831858407079595, 0.49576301487414098, 0.77345132743362799)
    6  = @(0.86581046211587998, 0.87429643527204504, 0.85925404644616399, 0.86491557223264504, 0.49387755102040798, 0.81613508442776705)
    7  = @(0.80497343013247502, 0.87018255578093295, 0.85199461118179698, 0.884381338742393,   0.49077539106354301, 0.87423935091277805)
    8  = @(0.88541232095675104, 0.90295358649789004, 0.88726649052055795, 0.92194092827004204, 0.49170259120196502, 0.892405063291139)
    9  = @(0.90392945544554404, 0.90990990990990905, 0.90798267326732596, 0.92342342342342298, 0.49900990099009901, 0.90990990990990905)
    10 = @(0.90530925013683605, 0.93857493857493801, 0.87657361795292799, 0.94594594594594505, 0.48617952928297697, 0.92874692874692799)
    11 = @(0.91582747482424398, 0.94256756756756699, 0.89891696750902494, 0.95270270270270196, 0.47814934448033403, 0.93581081081080997)
    12 = @(0.63468468468468398, 0.93401015228426398, 0.74684684684684599, 0.93908629441624303, 0.44999999999999901, 0.93908629441624303)
}

$newDataCols = @("C", "F", "I", "L", "O", "R")

foreach ($row in $newData.Keys) {
    $values = $newData[$row]
    for ($i = 0; $i -lt $newDataCols.Length; $i++) {
        $ws.Range("$($newDataCols[$i])$row").Value = $values[$i]
    }
}

# New "n" column: count of students used for the confidence interval at
# each semester.
$ws.Range("U2").Value = "n"

$nValues = @{
    3  = 689
    4  = 628
    5  = 565
    6  = 533
    7  = 493
    8  = 474
    9  = 444
    10 = 407
    11 = 296
    12 = 197
}

foreach ($row in $nValues.Keys) {
    $ws.Range("U$row").Value = $nValues[$row]
}

# Page orientation stays portrait (this also makes Excel emit the
# <pageSetup> element for the sheet).
$ws.PageSetup.Orientation = 1

$ws.Range("R16").Select() | Out-Null
